# Generate Report for Handoff
# Regenerates the localization-status report for a new handback run:
# the handed-back file's GUID changes, and the handoff/handback
# timestamps advance to the new run's values.

$wb = $excel.ActiveWorkbook

$newGuid = "46e19c79-0cb8-458e-83cf-68d8ba758c50"

$newZhXlf = "$newGuid.37636619180e1a395848851f9f141b10518df373.zh-cn.xlf"
$newDeXlf = "$newGuid.37636619180e1a395848851f9f141b10518df373.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-19 12:56:13"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newGuid.md"
}

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-19 12:56:07"
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-19 12:56:13"
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}
